$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 492. Excel shifts the
# existing rows 492:501 down to 495:504, carrying their values and
# formatting (including the date-formatted style in column D) along.
$ws.Rows("492:494").Insert()

# --- New row 492 ---
$ws.Cells.Item(492,1).Value = 5
$ws.Cells.Item(492,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(492,3).Value = "Maule"
$ws.Cells.Item(492,4).Value = 44628
$ws.Cells.Item(492,5).Value = 7
$ws.Cells.Item(492,6).Value = "Fruta"
$ws.Cells.Item(492,7).Value = 100101
$ws.Cells.Item(492,8).Value = "Berries"
$ws.Cells.Item(492,9).Value = 100112025
$ws.Cells.Item(492,10).Value = "Frutilla"
$ws.Cells.Item(492,11).Value = "Sin especificar"
$ws.Cells.Item(492,12).Value = "Especial"
$ws.Cells.Item(492,13).Value = 100
$ws.Cells.Item(492,14).Value = 8000
$ws.Cells.Item(492,15).Value = 8000
$ws.Cells.Item(492,16).Value = 8000
$ws.Cells.Item(492,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(492,18).Value = "Región del Maule"
$ws.Cells.Item(492,19).Value = 1143
$ws.Cells.Item(492,20).Value = 7

# --- New row 493 ---
$ws.Cells.Item(493,1).Value = 5
$ws.Cells.Item(493,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(493,3).Value = "Maule"
$ws.Cells.Item(493,4).Value = 44628
$ws.Cells.Item(493,5).Value = 7
$ws.Cells.Item(493,6).Value = "Fruta"
$ws.Cells.Item(493,7).Value = 100101
$ws.Cells.Item(493,8).Value = "Berries"
$ws.Cells.Item(493,9).Value = 100112025
$ws.Cells.Item(493,10).Value = "Frutilla"
$ws.Cells.Item(493,11).Value = "Sin especificar"
$ws.Cells.Item(493,12).Value = "Primera"
$ws.Cells.Item(493,13).Value = 120
$ws.Cells.Item(493,14).Value = 7000
$ws.Cells.Item(493,15).Value = 7000
$ws.Cells.Item(493,16).Value = 7000
$ws.Cells.Item(493,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(493,18).Value = "Región del Maule"
$ws.Cells.Item(493,19).Value = 1000
$ws.Cells.Item(493,20).Value = 7

# --- New row 494 ---
$ws.Cells.Item(494,1).Value = 5
$ws.Cells.Item(494,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(494,3).Value = "Maule"
$ws.Cells.Item(494,4).Value = 44628
$ws.Cells.Item(494,5).Value = 7
$ws.Cells.Item(494,6).Value = "Fruta"
$ws.Cells.Item(494,7).Value = 100101
$ws.Cells.Item(494,8).Value = "Berries"
$ws.Cells.Item(494,9).Value = 100112025
$ws.Cells.Item(494,10).Value = "Frutilla"
$ws.Cells.Item(494,11).Value = "Sin especificar"
$ws.Cells.Item(494,12).Value = "Primera"
$ws.Cells.Item(494,13).Value = 170
$ws.Cells.Item(494,14).Value = 6000
$ws.Cells.Item(494,15).Value = 6000
$ws.Cells.Item(494,16).Value = 6000
$ws.Cells.Item(494,17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(494,18).Value = "Región del Maule"
$ws.Cells.Item(494,19).Value = 857
$ws.Cells.Item(494,20).Value = 7
